$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value. D-column (price) values are forced to
# text so Excel does not reinterpret them as numbers/dates (they use "."
# as a separator and carry significant trailing zeros).
$updates = @(
    @('D2', '27.378.59'),
    @('E2', '  +2.57%  '),
    @('D3', '1.796.97'),
    @('E3', '  +3.54%  '),
    @('E4', '  +0.56%  '),
    @('D5', '337.28'),
    @('E5', '  +1.44%  '),
    @('D6', '1.0000'),
    @('E6', '  +0.35%  '),
    @('E7', '  +2.11%  '),
    @('D8', '0.3462'),
    @('E8', '  +1.83%  '),
    @('D9', '48.33'),
    @('E9', '  +0.32%  '),
    @('D10', '1.203'),
    @('E10', '  +1.34%  '),
    @('D11', '0.07511'),
    @('E11', '  +0.63%  '),
    @('D12', '1.000'),
    @('E12', '  +0.26%  '),
    @('D13', '22.02'),
    @('E13', '  +8.86%  '),
    @('D14', '6.490'),
    @('E14', '  +1.05%  '),
    @('D15', '1.795.68'),
    @('E15', '  +3.63%  '),
    @('D16', '7.052'),
    @('E16', '  -0.34%  '),
    @('D17', '0.00001097'),
    @('E17', '  +1.82%  '),
    @('D18', '0.06655'),
    @('E18', '  -0.98%  '),
    @('D19', '85.15'),
    @('E19', '  +3.18%  '),
    @('E20', '  +0.40%  '),
    @('D21', '6.517'),
    @('E21', '  +4.44%  '),
    @('D22', '17.35'),
    @('E22', '  +4.28%  '),
    @('D23', '27.374.92'),
    @('E23', '  +2.64%  '),
    @('D24', '12.52'),
    @('E24', '  -1.75%  '),
    @('D25', '2.436'),
    @('E25', '  -0.29%  '),
    @('D26', '2.573'),
    @('E26', '  +5.84%  '),
    @('D27', '1.496'),
    @('E27', '  -0.60%  '),
    @('D28', '21.39'),
    @('E28', '  +9.80%  '),
    @('D29', '152.82'),
    @('E29', '  +1.12%  '),
    @('D30', '1.996.90'),
    @('E30', '  +3.67%  '),
    @('D31', '134.20'),
    @('E31', '  +1.64%  '),
    @('D32', '4.057'),
    @('E32', '  -1.25%  '),
    @('D33', '6.138'),
    @('E33', '  +1.57%  '),
    @('D34', '0.08707'),
    @('E34', '  +1.31%  '),
    @('D35', '13.24'),
    @('E35', '  +3.19%  '),
    @('D36', '1.699'),
    @('E36', '  -0.02%  '),
    @('B37', 'InternetComputer(DFINITY)'),
    @('C37', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @('D37', '5.447'),
    @('E37', '  +0.81%  '),
    @('B38', 'TheSandbox'),
    @('C38', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('D38', '0.6900'),
    @('E38', '  +10.59%  '),
    @('D39', '8.947'),
    @('E39', '  +5.26%  '),
    @('B40', 'Algorand'),
    @('C40', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'),
    @('D40', '0.2213'),
    @('E40', '  +1.67%  '),
    @('B41', 'Hedera'),
    @('C41', 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @('D41', '0.06369'),
    @('E41', '  +2.24%  '),
    @('E42', '  -0.31%  '),
    @('E43', '  +4.07%  '),
    @('D44', '14.43'),
    @('E44', '  +0.89%  '),
    @('D45', '0.6465'),
    @('E45', '  +6.62%  '),
    @('D46', '0.9997'),
    @('E46', '  +0.37%  '),
    @('D47', '3.870'),
    @('E47', '  -0.88%  '),
    @('D48', '2.123'),
    @('E48', '  +3.36%  '),
    @('D49', '129.93'),
    @('E49', '  +0.65%  '),
    @('D50', '0.07194'),
    @('E50', '  -0.05%  '),
    @('D51', '79.42'),
    @('E51', '  +2.60%  ')
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $cell = $ws.Range($addr)
    if ($addr.Substring(0, 1) -eq "D") {
        $cell.NumberFormat = '@'
        $cell.Value = $val
        $cell.ClearFormats()
    } else {
        $cell.Value = $val
    }
}
